$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before C ("description"), shifting C:L -> D:M.
# Excel's Insert() carries the formatting of the column to the left (B)
# into the new blank column, which matches the target file's styling.
$ws.Columns("C").Insert()

# Header for the new "uom" column (old "mou" column F1, now shifted to F1)
# must be added to the shared-string table before "name_mm" so the
# sharedStrings order matches (uom, then name_mm at the tail).
$ws.Range("F1").Value = "uom"

# Header + data for the newly inserted "name_mm" column: duplicate the
# "name" column (B) values, matching the source edit (localized name not
# yet translated, so it mirrors the English name).
$ws.Range("C1").Value = "name_mm"
for ($r = 2; $r -le 21; $r++) {
  $nameValue = $ws.Cells.Item($r, 2).Value2
  $ws.Cells.Item($r, 3).Value = $nameValue
}

# Match the author's final selection (cell C2, the first data cell of the
# newly-added column) and let the view default back to the top-left.
[void]$ws.Range("C2").Select()
